$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 07:48:01"
$ws.Range("E3").Value = "2026-02-06 07:48:03"
$ws.Range("G3").Value = "174 cm"
$ws.Range("O3").Value = "-2.4 °C"
$ws.Range("E4").Value = "2026-02-06 07:48:06"
$ws.Range("J4").Value = "994.2 hPa"
$ws.Range("O4").Value = "11.4 °C"
$ws.Range("E5").Value = "2026-02-06 07:48:09"
$ws.Range("J5").Value = "994.6 hPa"
$ws.Range("O5").Value = "7.7 °C"
$ws.Range("E6").Value = "2026-02-06 07:48:11"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "52%"
$ws.Range("J6").Value = "995.7 hPa"
$ws.Range("K6").Value = "0.0 MJ/m2"
$ws.Range("N6").Value = "13.4 °C 7:00 TU"
$ws.Range("O6").Value = "14.1 °C"
$ws.Range("E7").Value = "2026-02-06 07:48:13"
$ws.Range("J7").Value = "995.5 hPa"
$ws.Range("O7").Value = "9.7 °C"
$ws.Range("E8").Value = "2026-02-06 07:48:16"
$ws.Range("K8").Value = "0.0 MJ/m2"
$ws.Range("N8").Value = "3.9 °C 7:15 TU"
$ws.Range("E9").Value = "2026-02-06 07:48:19"
$ws.Range("N9").Value = "-0.1 °C 7:25 TU"
$ws.Range("O9").Value = "1.6 °C"
$ws.Range("E10").Value = "2026-02-06 07:48:21"
$ws.Range("E11").Value = "2026-02-06 07:48:23"
$ws.Range("J11").Value = "996.7 hPa"
$ws.Range("N11").Value = "0.0 °C 7:17 TU"
$ws.Range("O11").Value = "3.5 °C"
$ws.Range("E12").Value = "2026-02-06 07:48:26"
$ws.Range("O12").Value = "11.9 °C"
$ws.Range("E13").Value = "2026-02-06 07:48:29"
$ws.Range("O13").Value = "5.6 °C"
$ws.Range("E14").Value = "2026-02-06 07:48:31"
$ws.Range("I14").Value = "0.4 mm"
$ws.Range("N14").Value = "-5.2 °C 7:25 TU"
$ws.Range("E15").Value = "2026-02-06 07:48:34"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "90%"
$ws.Range("J15").Value = "994.8 hPa"
$ws.Range("O15").Value = "5.4 °C"
$ws.Range("E16").Value = "2026-02-06 07:48:37"
$ws.Range("L16").Value = "13.3 km/h - 324º 7:09 TU"
$ws.Range("E17").Value = "2026-02-06 07:48:39"
$ws.Range("J17").Value = "997.6 hPa"
$ws.Range("O17").Value = "2.7 °C"
$ws.Range("E18").Value = "2026-02-06 07:48:42"
$ws.Range("N18").Value = "-5.5 °C 7:29 TU"
$ws.Range("E19").Value = "2026-02-06 07:48:45"
$ws.Range("J19").Value = "997.9 hPa"
$ws.Range("E20").Value = "2026-02-06 07:48:47"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "72%"
$ws.Range("K20").Value = "0.0 MJ/m2"
$ws.Range("N20").Value = "-5.0 °C 7:03 TU"
$ws.Range("O20").Value = "-2.6 °C"
$ws.Range("E21").Value = "2026-02-06 07:48:50"
$ws.Range("J21").Value = "995.7 hPa"
$ws.Range("K21").Value = "0.0 MJ/m2"
$ws.Range("N21").Value = "1.7 °C 7:04 TU"
$ws.Range("O21").Value = "4.1 °C"
$ws.Range("E22").Value = "2026-02-06 07:48:53"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "87%"
$ws.Range("K22").Value = "0.0 MJ/m2"
$ws.Range("N22").Value = "3.2 °C 7:05 TU"
$ws.Range("O22").Value = "6.7 °C"
$ws.Range("E23").Value = "2026-02-06 07:48:55"
$ws.Range("J23").Value = "994.8 hPa"
$ws.Range("E24").Value = "2026-02-06 07:48:58"
$ws.Range("J24").Value = "993.7 hPa"
$ws.Range("K24").Value = "0.0 MJ/m2"
$ws.Range("E25").Value = "2026-02-06 07:49:01"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "92%"
$ws.Range("J25").Value = "996.9 hPa"
$ws.Range("E26").Value = "2026-02-06 07:49:04"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "81%"
$ws.Range("O26").Value = "-1.5 °C"
$ws.Range("E27").Value = "2026-02-06 07:49:07"
$ws.Range("J27").Value = "994.5 hPa"
$ws.Range("N27").Value = "4.6 °C 7:07 TU"
$ws.Range("O27").Value = "6.9 °C"
$ws.Range("E28").Value = "2026-02-06 07:49:09"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "92%"
$ws.Range("J28").Value = "998.0 hPa"
$ws.Range("O28").Value = "1.7 °C"
$ws.Range("E29").Value = "2026-02-06 07:49:12"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "70%"
$ws.Range("K29").Value = "0.0 MJ/m2"
$ws.Range("N29").Value = "5.3 °C 7:14 TU"
$ws.Range("O29").Value = "9.8 °C"
$ws.Range("E30").Value = "2026-02-06 07:49:15"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "71%"
$ws.Range("K30").Value = "0.0 MJ/m2"
$ws.Range("O30").Value = "-3.9 °C"
$ws.Range("E31").Value = "2026-02-06 07:49:17"
$ws.Range("J31").Value = "997.5 hPa"
$ws.Range("E32").Value = "2026-02-06 07:49:19"
$ws.Range("J32").Value = "996.0 hPa"
$ws.Range("E33").Value = "2026-02-06 07:49:22"
$ws.Range("O33").Value = "6.1 °C"
$ws.Range("E34").Value = "2026-02-06 07:49:25"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "86%"
$ws.Range("K34").Value = "0.0 MJ/m2"
$ws.Range("N34").Value = "2.2 °C 7:02 TU"
$ws.Range("O34").Value = "6.3 °C"
$ws.Range("E35").Value = "2026-02-06 07:49:27"
$ws.Range("N35").Value = "-3.4 °C 7:08 TU"
$ws.Range("E36").Value = "2026-02-06 07:49:30"
$ws.Range("J36").Value = "997.5 hPa"
$ws.Range("O36").Value = "10.6 °C"
